$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Game_Record")
$ws.Select()
$ws.Range("A13").Formula = "=ROW()-1"
$ws.Range("B13").Value = 45923
$ws.Range("C13").Value = "Doanage"
$ws.Range("D13").Value = "Player1"
$ws.Range("E13").Value = "SimpleJack"
$ws.Range("F13").Value = "DrSystomatix"
$ws.Range("C14").Select()
